# Updated symbol list on Mon Dec 26 08:47:04 UTC 2022 with GitHub Actions
#
# Refresh the crypto price ticker sheet:
#  - bump the "Price" (column D) quote for most rows to the latest pull
#  - "One" (ONE) jumped to the top of the TigerCash..BTSEToken block
#    (rows 18-24), pushing the rest of that block down by one row, and
#    the "Best in 24h" tag moved along with the reshuffle
#
# Note: column D stores plain text (e.g. "243.10", "0.03300") rather than
# numbers, so values are entered with a leading apostrophe to force text
# and then the cell style is reset to "Normal" so Excel doesn't leave the
# cell tagged with a '@' (Text) number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

# Rows 2-17: price-only refresh
Set-TextValue "D2"  "243.10"
Set-TextValue "D3"  "23.08"
Set-TextValue "D4"  "5.417"
Set-TextValue "D5"  "0.05919"
Set-TextValue "D7"  "6.547"
Set-TextValue "D8"  "0.8133"
Set-TextValue "D9"  "0.9103"
Set-TextValue "D10" "0.1405"
Set-TextValue "D11" "0.07417"
Set-TextValue "D12" "0.03300"
Set-TextValue "D13" "0.03052"
Set-TextValue "D14" "0.09353"
Set-TextValue "D15" "3.874"
Set-TextValue "D16" "0.001557"
Set-TextValue "D17" "0.04674"

# Row 18: TigerCash -> One (moved up from row 24)
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005943"
$ws.Range("E18").Value = "17OneONE"

# Row 19: HotbitToken -> TigerCash
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D19" "0.006104"
$ws.Range("E19").Value = "18TigerCashTCH"

# Row 20: BitKan -> HotbitToken (now the 24h best performer)
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.004980"
$ws.Range("E20").Value = "19HotbitTokenHTBBestin24h"

# Row 21: NitroEx -> BitKan
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.0009820"
$ws.Range("E21").Value = "20BitKanKAN"

# Row 22: LEO -> NitroEx
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.00009003"
$ws.Range("E22").Value = "21NitroExNTX"

# Row 23: BTSEToken -> LEO
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "3.605"
$ws.Range("E23").Value = "22LEOLEO"

# Row 24: One -> BTSEToken
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.136"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# Row 25 and remaining scattered price-only refreshes
Set-TextValue "D25" "0.3239"
Set-TextValue "D40" "0.04008"
Set-TextValue "D41" "0.006210"
Set-TextValue "D42" "0.1076"
Set-TextValue "D43" "0.002621"
Set-TextValue "D44" "0.008213"
Set-TextValue "D45" "0.00005243"
Set-TextValue "D49" "0.002259"
